# Generate Report for Handoff
# Update the "b.md" rows across the Overview / zh-cn / de-de sheets to
# reflect that the handoff has occurred ("Ready for handoff") with a new
# handoff file name and timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is "b.md" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-24-18 08:24:51"

# --- zh-cn sheet: row 3 is "b.md" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-18 08:24:49"

# --- de-de sheet: row 3 is "b.md" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-18 08:24:51"
